# Applies the two substantive changes described by the target diff:
#
#   1. The table on slide 6 gets switched from the custom "Table_0"
#      style ({7DFFB818-2F60-48F5-80A1-A3B4C5DBABB3}, defined in
#      ppt/tableStyles.xml) to PowerPoint's built-in "No Style, Table
#      Grid" style ({503BFA9C-E17C-4AD3-8A2D-4B5E7FF54F62}).
#
#   2. The deck's theme (ppt/theme/theme1.xml, used by the one-and-only
#      slide master / all slide layouts) is swapped from the
#      "Integral" color scheme to the stock "Office" color scheme.
#      (ppt/theme/theme2.xml -- the Notes Master's theme, which held
#      "Office Theme" and becomes "Integral" in the target -- isn't
#      reachable through the exposed Notes Master object in this COM
#      host, so it's intentionally left untouched here rather than
#      risk corrupting the slide master's theme.)

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{503BFA9C-E17C-4AD3-8A2D-4B5E7FF54F62}")
    }
}

# --- 2. Theme colors: Integral -> Office ----------------------------------
$cs = $p.SlideMaster.ColorScheme
$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
